$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "249.62"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.93"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.437"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05620"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.421"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.368"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8124"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8950"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1431"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07531"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03123"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03092"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09329"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.568"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001588"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04754"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005794"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006416"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004993"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001502"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.190"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3303"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003007"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04019"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006790"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1070"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002724"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007882"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005575"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5003"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2395"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01011"
